$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new shared strings / headers for GN1 (median) and GO1 (mean),
# copying the header style (bold, centered, bordered) from GM1.
$ws.Range("GM1").Copy()
$ws.Range("GN1:GO1").PasteSpecial(-4122)
$ws.Range("GN1").Value = "median"
$ws.Range("GO1").Value = "mean"

# Fill in per-year median/mean difference values (GN = median, GO = mean)
$ws.Range("GN2").Value = 0.06337205386256721
$ws.Range("GO2").Value = 0.06528782122633213
$ws.Range("GN3").Value = 0.06800478836873336
$ws.Range("GO3").Value = 0.08816384318337991
$ws.Range("GN4").Value = 0.04501940460501911
$ws.Range("GO4").Value = 0.05538780466068421
$ws.Range("GN5").Value = 0.05957911306220465
$ws.Range("GO5").Value = 0.06186230780697061
$ws.Range("GN6").Value = 0.03544842760778377
$ws.Range("GO6").Value = 0.04631255896461434
$ws.Range("GN7").Value = 0.05577216945653447
$ws.Range("GO7").Value = 0.05611702731305356
$ws.Range("GN8").Value = 0.05552195128335391
$ws.Range("GO8").Value = 0.06041647746370441
$ws.Range("GN9").Value = 0.05496370016980578
$ws.Range("GO9").Value = 0.06513491915690539
$ws.Range("GN10").Value = 0.08464078528023361
$ws.Range("GO10").Value = 0.09747536574614891
$ws.Range("GN11").Value = 0.06883488161428845
$ws.Range("GO11").Value = 0.07442556814735798
$ws.Range("GN12").Value = 0.05109069629356972
$ws.Range("GO12").Value = 0.07334332841332483
$ws.Range("GN13").Value = 0.04766735396282341
$ws.Range("GO13").Value = 0.06134053769369419
$ws.Range("GN14").Value = 0.0594874255599267
$ws.Range("GO14").Value = 0.07052519395066995
$ws.Range("GN15").Value = 0.08378539090723093
$ws.Range("GO15").Value = 0.09902487179227634
$ws.Range("GN16").Value = 0.04732755504470741
$ws.Range("GO16").Value = 0.05901726820832719
$ws.Range("GN17").Value = 0.08359735415185876
$ws.Range("GO17").Value = 0.08522391949284648
$ws.Range("GN18").Value = 0.04903142856846278
$ws.Range("GO18").Value = 0.05669827749372373
$ws.Range("GN19").Value = 0.05884624391667165
$ws.Range("GO19").Value = 0.0665494957153145
$ws.Range("GN20").Value = 0.08372931659904501
$ws.Range("GO20").Value = 0.08734079008810451
$ws.Range("GN21").Value = 0.08908607061578332
$ws.Range("GO21").Value = 0.08913828601802001
$ws.Range("GN22").Value = 0.08573158600711751
$ws.Range("GO22").Value = 0.09645674845840235
$ws.Range("GN23").Value = 0.07978618548158908
$ws.Range("GO23").Value = 0.08580643812549808
$ws.Range("GN24").Value = 0.09645140458355916
$ws.Range("GO24").Value = 0.09416952448753312
$ws.Range("GN25").Value = 0.09176422803523557
$ws.Range("GO25").Value = 0.09671711837747356
$ws.Range("GN26").Value = 0.08537684954474048
$ws.Range("GO26").Value = 0.09950345812318188
$ws.Range("GN27").Value = 0.09716210409719085
$ws.Range("GO27").Value = 0.1074028124585002
$ws.Range("GN28").Value = 0.08178861771441469
$ws.Range("GO28").Value = 0.09811507834023893
$ws.Range("GN29").Value = 0.1029815994457294
$ws.Range("GO29").Value = 0.1141841859551018
$ws.Range("GN30").Value = 0.09355639961441045
$ws.Range("GO30").Value = 0.1149904916624163
$ws.Range("GN31").Value = 0.09579196242741533
$ws.Range("GO31").Value = 0.1150678376295326
$ws.Range("GN32").Value = 0.1177236569173039
$ws.Range("GO32").Value = 0.1260044950220779
$ws.Range("GN33").Value = 0.1178626964275973
$ws.Range("GO33").Value = 0.1382125169329668
$ws.Range("GN34").Value = 0.1208353403092671
$ws.Range("GO34").Value = 0.1364601414520136
$ws.Range("GN35").Value = 0.1211903627676592
$ws.Range("GO35").Value = 0.1438918803800114
$ws.Range("GN36").Value = 0.1279560596765174
$ws.Range("GO36").Value = 0.154149618360005
$ws.Range("GN37").Value = 0.1302768565108977
$ws.Range("GO37").Value = 0.1590271494657305
$ws.Range("GN38").Value = 0.1323621667220299
$ws.Range("GO38").Value = 0.166424386785574
$ws.Range("GN39").Value = 0.1351326435847696
$ws.Range("GO39").Value = 0.1751967931967772
